$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows 2-49 ---
# Column D holds text-like price strings (e.g. '24.172.35', '0.9998').
# Force Text format before writing so Excel doesn't coerce numeric-looking
# strings (like '0.9998') into real numbers - matches the original inlineStr data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.172.35'
$ws.Range("E2").Value = '  -2.67%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.643.27'
$ws.Range("E3").Value = '  -2.60%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.56%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.04'
$ws.Range("E5").Value = '  -2.06%  '

$ws.Range("E6").Value = '  -0.49%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3903'
$ws.Range("E7").Value = '  -0.72%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3866'
$ws.Range("E8").Value = '  -2.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9999'
$ws.Range("E9").Value = '  -0.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.76'
$ws.Range("E10").Value = '  -3.96%  '

$ws.Range("E11").Value = '  -4.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08681'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.68'
$ws.Range("E13").Value = '  -6.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.121'
$ws.Range("E14").Value = '  -2.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001292'
$ws.Range("E15").Value = '  -2.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.458'
$ws.Range("E16").Value = '  -4.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.662.48'
$ws.Range("E17").Value = '  +4.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.96'
$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06898'
$ws.Range("E19").Value = '  -3.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.53'
$ws.Range("E20").Value = '  +1.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.898'
$ws.Range("E21").Value = '  -3.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.59'
$ws.Range("E23").Value = '  -3.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.165.71'
$ws.Range("E24").Value = '  -2.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.331'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.785'
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.39'
$ws.Range("E27").Value = '  -4.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.69'
$ws.Range("E28").Value = '  -2.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.480'
$ws.Range("E29").Value = '  +7.85%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '140.39'
$ws.Range("E30").Value = '  -6.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.358'
$ws.Range("E31").Value = '  -8.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.412'
$ws.Range("E32").Value = '  -8.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.776.13'
$ws.Range("E33").Value = '  -4.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.972'
$ws.Range("E34").Value = '  +0.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08067'
$ws.Range("E35").Value = '  -4.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02901'
$ws.Range("E36").Value = '  -5.82%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2681'
$ws.Range("E37").Value = '  -4.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9495'
$ws.Range("E38").Value = '  -6.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09215'
$ws.Range("E39").Value = '  -3.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.04'
$ws.Range("E40").Value = '  -5.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.457'
$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7545'
$ws.Range("E42").Value = '  -5.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.02'
$ws.Range("E43").Value = '  -5.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.98'
$ws.Range("E44").Value = '  -4.63%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6905'
$ws.Range("E45").Value = '  -3.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.463'
$ws.Range("E46").Value = '  -4.96%  '

$ws.Range("E47").Value = '  -2.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9986'
$ws.Range("E48").Value = '  -0.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08411'
$ws.Range("E49").Value = '  -3.60%  '

# --- Rows 50 and 51 swapped positions (Flow <-> Quant) with updated data ---
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.85'
$ws.Range("E50").Value = '  -3.62%  '

$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.262'
$ws.Range("E51").Value = '  -5.92%  '
